$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.870.00'
$ws.Range('E2').Value = '  -1.65%  '
$ws.Range('D3').Value = '1.807.80'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.51%  '
$ws.Range('D5').Value = '''309.94'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.40%  '
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').Value = '''0.4444'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +4.60%  '
$ws.Range('D8').Value = '''0.3668'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.33%  '
$ws.Range('D9').Value = '''0.07329'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.87%  '
$ws.Range('D10').Value = '''0.8553'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.34%  '
$ws.Range('D11').Value = '''20.68'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.23%  '
$ws.Range('D12').Value = '1.814.69'
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('D13').Value = '''6.603'
$ws.Range('D13').Style = "Normal"
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').Value = '''92.22'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.68%  '
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D15').Value = '''0.07088'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('D16').Value = '''5.311'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').Value = '''0.000008716'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.89%  '
$ws.Range('E19').Value = '  -0.41%  '
$ws.Range('D20').Value = '''14.88'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.54%  '
$ws.Range('D21').Value = '26.883.89'
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('D22').Value = '''5.158'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.44%  '
$ws.Range('D23').Value = '''10.86'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.53%  '
$ws.Range('D24').Value = '''1.989'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('D25').Value = '''151.60'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('D26').Value = '''2.208'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.18%  '
$ws.Range('D27').Value = '''18.63'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.10%  '
$ws.Range('D28').Value = '''5.200'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('D29').Value = '''116.85'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').Value = '''0.08820'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.56%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''0.7523'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = '''1.176'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.09%  '
$ws.Range('D33').Value = '''2.931'
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Value = '''4.455'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('D35').Value = '''0.9997'
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Value = '''1.091'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.63%  '
$ws.Range('D37').Value = '''0.01968'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('D38').Value = '''0.05195'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('D39').Value = '''0.5365'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +5.82%  '
$ws.Range('D40').Value = '''2.865'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.63%  '
$ws.Range('D41').Value = '''7.042'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.28%  '
$ws.Range('D42').Value = '''0.1688'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('D43').Value = '''0.5230'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +9.61%  '
$ws.Range('D44').Value = '''8.443'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.99%  '
$ws.Range('D45').Value = '''10.60'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.36%  '
$ws.Range('D46').Value = '''1.970'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +5.43%  '
$ws.Range('D47').Value = '''105.46'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.15%  '
$ws.Range('D48').Value = '''0.9995'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('D49').Value = '''1.670'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('D50').Value = '''0.06337'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').Value = '''0.9208'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.33%  '
